$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 7943
$ws1.Range("F6").Value = 6
$ws1.Range("F8").Value = 1897
$ws1.Range("F14").Value = 1063
$ws1.Range("F23").Value = 584
$ws1.Range("F24").Value = 1205
$ws1.Range("F25").Value = 1081
$ws1.Range("F26").Value = 609
$ws1.Range("F31").Value = 124
$ws1.Range("F41").Value = 531
$ws1.Range("F44").Value = 727
$ws1.Range("F45").Value = 65

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 356
$ws2.Range("F16").Value = 11

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value = 1419
$ws3.Range("F7").Value = 603
$ws3.Range("F8").Value = 2260
$ws3.Range("F9").Value = 9090

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 1419
$ws4.Range("F6").Value = 2260
$ws4.Range("F14").Value = 1063
$ws4.Range("F21").Value = 584
$ws4.Range("F22").Value = 609
$ws4.Range("F28").Value = 11
$ws4.Range("F38").Value = 531
$ws4.Range("F41").Value = 727
$ws4.Range("F43").Value = 65
